$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 700.4783
$ws.Range("I19").Value = 379.85715
$ws.Range("J19").Value = 840.75
$ws.Range("K19").Value = 379.85715
$ws.Range("L19").Value = 840.75
$ws.Range("M19").Value = -204.85715
$ws.Range("N19").Value = -1190.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 9966.75
$ws.Range("I113").Value = 7870.1055
$ws.Range("J113").Value = 13031.077
$ws.Range("K113").Value = 7870.1055
$ws.Range("L113").Value = 13031.077
$ws.Range("M113").Value = -4616.1055
$ws.Range("N113").Value = -19539.077

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H134").Value = 125000
$ws.Range("J134").Value = 125000
$ws.Range("L134").Value = 125000
$ws.Range("N134").Value = -135140

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1124.1072
$ws.Range("I135").Value = 908.03845
$ws.Range("J135").Value = 3933
$ws.Range("K135").Value = 8172.34605
$ws.Range("L135").Value = 35397
$ws.Range("M135").Value = -5637.34605
$ws.Range("N135").Value = -40467

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H136").Value = 150000
$ws.Range("J136").Value = 150000
$ws.Range("L136").Value = 150000
$ws.Range("N136").Value = -160200

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1912.7826
$ws.Range("I137").Value = 1871.2858
$ws.Range("J137").Value = 1977.3334
$ws.Range("K137").Value = 5613.857400000001
$ws.Range("L137").Value = 5932.0002
$ws.Range("M137").Value = -3063.857400000001
$ws.Range("N137").Value = -11032.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10233.071
$ws.Range("I32").Value = 4894.49
$ws.Range("J32").Value = 24562.947
$ws.Range("K32").Value = 4894.49
$ws.Range("L32").Value = 24562.947
$ws.Range("M32").Value = -4607.49
$ws.Range("N32").Value = -25136.947

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 4224.5
$ws.Range("I102").Value = 4900
$ws.Range("K102").Value = 4900
$ws.Range("M102").Value = -3278

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H119").Value = 100000
$ws.Range("J119").Value = 100000
$ws.Range("L119").Value = 100000
$ws.Range("N119").Value = -109676

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1961.8125
$ws.Range("I122").Value = 1728.8462
$ws.Range("K122").Value = 5186.5386
$ws.Range("M122").Value = -2736.5386

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 16670

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6332.8237
$ws.Range("I86").Value = 6077.615
$ws.Range("K86").Value = 6077.615
$ws.Range("M86").Value = -4954.615

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 6332.8237
$ws.Range("I89").Value = 6077.615
$ws.Range("K89").Value = 30388.075
$ws.Range("M89").Value = -24772.075

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1361.25
$ws.Range("I107").Value = 1138.3125
$ws.Range("K107").Value = 1138.3125
$ws.Range("M107").Value = 781.6875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 69217.13
$ws.Range("I31").Value = 85274
$ws.Range("J31").Value = 4989.6665
$ws.Range("K31").Value = 85274
$ws.Range("L31").Value = 4989.6665
$ws.Range("M31").Value = -84979
$ws.Range("N31").Value = -5579.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 69217.13
$ws.Range("I34").Value = 85274
$ws.Range("J34").Value = 4989.6665
$ws.Range("K34").Value = 85274
$ws.Range("L34").Value = 4989.6665
$ws.Range("M34").Value = -85072
$ws.Range("N34").Value = -5393.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2575.85
$ws.Range("I58").Value = 1508.6666
$ws.Range("K58").Value = 1508.6666
$ws.Range("M58").Value = -1305.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2311.875
$ws.Range("I122").Value = 1518.3
$ws.Range("J122").Value = 3634.5
$ws.Range("K122").Value = 4554.9
$ws.Range("L122").Value = 10903.5
$ws.Range("M122").Value = -2104.9
$ws.Range("N122").Value = -15803.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H125").Value = 56665
$ws.Range("J125").Value = 56665
$ws.Range("L125").Value = 56665
$ws.Range("N125").Value = -61585

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H127").Value = 100780
$ws.Range("J127").Value = 100780
$ws.Range("L127").Value = 100780
$ws.Range("N127").Value = -110700

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3931.889
$ws.Range("I132").Value = 3886.0435
$ws.Range("K132").Value = 11658.1305
$ws.Range("M132").Value = -9128.130500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2575.85
$ws.Range("I136").Value = 1508.6666
$ws.Range("K136").Value = 4525.9998
$ws.Range("M136").Value = -1975.9998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 587581.5
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 12333.333
$ws.Range("I110").Value = 12333.333
$ws.Range("K110").Value = 36999.999
$ws.Range("M110").Value = -32909.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 6250771
$ws.Range("I113").Value = 8334046.5
$ws.Range("J113").Value = 945
$ws.Range("K113").Value = 25002139.5
$ws.Range("L113").Value = 2835
$ws.Range("M113").Value = -24999969.5
$ws.Range("N113").Value = -7175

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2453.0715
$ws.Range("I80").Value = 2259.7778
$ws.Range("K80").Value = 2259.7778
$ws.Range("M80").Value = -1261.7778

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2453.0715
$ws.Range("I83").Value = 2259.7778
$ws.Range("K83").Value = 11298.889
$ws.Range("M83").Value = -6306.888999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 500000540
$ws.Range("I102").Value = 1100
$ws.Range("K102").Value = 1100
$ws.Range("M102").Value = 522

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4231.8887
$ws.Range("I113").Value = 4440.875
$ws.Range("J113").Value = 2560
$ws.Range("K113").Value = 4440.875
$ws.Range("L113").Value = 2560
$ws.Range("M113").Value = -2270.875
$ws.Range("N113").Value = -6900

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H121").Value = 19999.334
$ws.Range("J121").Value = 19999.334
$ws.Range("L121").Value = 19999.334
$ws.Range("N121").Value = -23493.334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4989.615
$ws.Range("I122").Value = 5238.6
$ws.Range("J122").Value = 4834
$ws.Range("K122").Value = 15715.8
$ws.Range("L122").Value = 14502
$ws.Range("M122").Value = -13265.8
$ws.Range("N122").Value = -19402

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H70").Value = 25000
$ws.Range("J70").Value = 25000
$ws.Range("L70").Value = 25000
$ws.Range("N70").Value = -25540

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H73").Value = 25000
$ws.Range("J73").Value = 25000
$ws.Range("L73").Value = 25000
$ws.Range("N73").Value = -26872

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4491.278
$ws.Range("I122").Value = 4178.5835
$ws.Range("K122").Value = 12535.7505
$ws.Range("M122").Value = -10085.7505

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1699.9032
$ws.Range("I122").Value = 1636.32
$ws.Range("K122").Value = 4908.96
$ws.Range("M122").Value = -2458.96

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
